$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" (Changed) date column C was bumped by one day
# (serial 45177 -> 45178, i.e. 2023-09-08 -> 2023-09-09) for every
# data row (rows 2 through 51).
$ws.Range("C2:C51").Value2 = 45178
